# no-op test
$wb = $excel.ActiveWorkbook
Write-Host "sheets:" $wb.Worksheets.Count
